$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.812.51'
$ws.Range('E2').Value = '  -1.16%  '
$ws.Range('D3').Value = '3.495.65'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '602.44'
$ws.Range('E5').Value = '  -1.11%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '147.95'
$ws.Range('E6').Value = '  -2.45%  '
$ws.Range('D7').Value = '3.494.66'
$ws.Range('E7').Value = '  -0.40%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.479'
$ws.Range('E9').Value = '  -1.52%  '
$ws.Range('E10').Value = '  -0.61%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '7.89'
$ws.Range('E11').Value = '  +3.49%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.422'
$ws.Range('E12').Value = '  -2.68%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000213'
$ws.Range('E13').Value = '  -1.47%  '
$ws.Range('D14').Value = '4.084.92'
$ws.Range('E14').Value = '  -0.38%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '31.18'
$ws.Range('E15').Value = '  -4.10%  '
$ws.Range('D16').Value = '3.485.91'
$ws.Range('E16').Value = '  -0.67%  '
$ws.Range('D17').Value = '66.755.34'
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.116'
$ws.Range('E18').Value = '  -0.80%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '10.47'
$ws.Range('E19').Value = '  +6.16%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.36'
$ws.Range('E20').Value = '  -3.00%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '15.30'
$ws.Range('E21').Value = '  -1.87%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '433.52'
$ws.Range('E22').Value = '  -3.34%  '
$ws.Range('E23').Value = '  -4.31%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '79.69'
$ws.Range('E24').Value = '  +2.02%  '
$ws.Range('D25').Value = '3.632.06'
$ws.Range('E25').Value = '  -0.32%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('E27').Value = '  -3.07%  '
$ws.Range('E28').Value = '  -6.88%  '
$ws.Range('E29').Value = '  -2.86%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '8.21'
$ws.Range('E30').Value = '  -7.32%  '
$ws.Range('E31').Value = '  -0.89%  '
$ws.Range('E32').Value = '  -2.86%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.998'
$ws.Range('E33').Value = '  -0.20%  '
$ws.Range('E34').Value = '  -2.01%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '25.33'
$ws.Range('E35').Value = '  -1.66%  '
$ws.Range('D36').Value = '3.486.71'
$ws.Range('E36').Value = '  -0.41%  '
$ws.Range('E37').Value = '  -5.04%  '
$ws.Range('E38').Value = '  -4.92%  '
$ws.Range('B39').Value = 'Aptos'
$ws.Range('C39').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '7.98'
$ws.Range('E39').Value = '  -0.69%  '
$ws.Range('B40').Value = 'USDe'
$ws.Range('C40').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.00'
$ws.Range('E40').Value = '  +0.03%  '
$ws.Range('E41').Value = '  -0.05%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.0891'
$ws.Range('E42').Value = '  -0.81%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '169.86'
$ws.Range('E43').Value = '  -2.39%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.10'
$ws.Range('E44').Value = '  -9.45%  '
$ws.Range('E45').Value = '  -1.34%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.896'
$ws.Range('E46').Value = '  +1.73%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '45.67'
$ws.Range('E47').Value = '  -1.99%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '28.39'
$ws.Range('E48').Value = '  -6.44%  '
$ws.Range('E49').Value = '  +0.22%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '7.45'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '2.42'
$ws.Range('E51').Value = '  -3.78%  '
